# Alternate mapping for panel2 (panel2_v2 sheet): insert an alternate
# FITC-A/BV711-A based gating pair ahead of the existing BB515-A/BV711-A
# based "Myeloid DC" / "Plasmacytoid DC" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("panel2_v2")
$ws.Activate()

# Make room for two new rows right before the old row 13 (shifts the
# three trailing NK rows down from 13-15 to 15-17, matching the diff).
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# New row 13: FITC-A+BV711-A- -> Dendritic -> Myeloid DC (CD11c+ CD123-)
$ws.Range("A13").Value = "FITC-A+BV711-A-"
$ws.Range("B13").Value = "Dendritic"
$ws.Range("C13").Value = "Myeloid DC (CD11c+ CD123-)"

# New row 14: FITC-A-BV711-A- -> Dendritic -> Plasmacytoid DC (CD11c- CD123+)
$ws.Range("A14").Value = "FITC-A-BV711-A-"
$ws.Range("B14").Value = "Dendritic"
$ws.Range("C14").Value = "Plasmacytoid DC (CD11c- CD123+)"

# Match the saved selection/scroll state recorded in the workbook.
$excel.Goto($ws.Range("A4"), $true)
$ws.Range("C9").Select()
